$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1) Re-point cell formatting BEFORE overwriting values, by copying
#    existing formats around (this reuses existing style records
#    instead of synthesizing new ones with slightly different XML).
# -----------------------------------------------------------------

# E2 currently carries the "code" style (bold-ish black font + thin
# border) that the new layout needs on columns C (Emp ID) and G
# (Client). Clone it there first, while E2 still has it.
$ws.Range("E2").Copy() | Out-Null
$ws.Range("C2:C3").PasteSpecial(-4122) | Out-Null
$ws.Range("E2").Copy() | Out-Null
$ws.Range("G2:G3").PasteSpecial(-4122) | Out-Null
$ws.Range("E2").Copy() | Out-Null
$ws.Range("D2:D3").PasteSpecial(-4122) | Out-Null

# New column D (Assignee_QA) keeps that same font but loses its left
# edge (it sits right next to column C which already has a right/top/
# bottom/left boxed border).
$ws.Range("D2:D3").Borders.Item(7).LineStyle = -4142

# E (Typist) and F (Typist QC) go back to the plain data style.
$ws.Range("F2").Copy() | Out-Null
$ws.Range("E2:E3").PasteSpecial(-4122) | Out-Null

# New column N (Tier): header uses the shared header style, the two
# data cells use the plain data style (and stay empty).
$ws.Range("M1").Copy() | Out-Null
$ws.Range("N1").PasteSpecial(-4122) | Out-Null
$ws.Range("M2").Copy() | Out-Null
$ws.Range("N2:N3").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# -----------------------------------------------------------------
# 2) Header row (row 1): columns got reshuffled + "Tier" appended.
# -----------------------------------------------------------------
$ws.Range("A1").Value = "Order Received Data and Time"
$ws.Range("B1").Value = "OrderID"
$ws.Range("C1").Value = "Emp ID-Order Assigned"
$ws.Range("D1").Value = "Assignee_QA"
$ws.Range("E1").Value = "Typist"
$ws.Range("F1").Value = "Typist QC"
$ws.Range("G1").Value = "Client"
$ws.Range("H1").Value = "Lob"
$ws.Range("I1").Value = "Process"
$ws.Range("J1").Value = "Product Name"
$ws.Range("K1").Value = "State"
$ws.Range("L1").Value = "County"
$ws.Range("M1").Value = "Status"
$ws.Range("N1").Value = "Tier"

# -----------------------------------------------------------------
# 3) Row 2 (new order FINN18-001)
# -----------------------------------------------------------------
$ws.Range("A2").Value = 45436.041666666664
$ws.Range("B2").Value = "FINN18-001"
$ws.Range("C2").Value = "SIPL4167"
$ws.Range("D2").Value = "SIPL5688"
$ws.Range("E2").Value = "SIPL5317"
$ws.Range("F2").Value = "SIPL5317"
$ws.Range("G2").Value = "FINN TITLE"
$ws.Range("H2").Value = "Title"
$ws.Range("I2").Value = "Search & Typing"
$ws.Range("J2").Value = "Property Reports"
$ws.Range("K2").Value = "AL"
$ws.Range("L2").Value = "Shelby"
$ws.Range("M2").Value = "WIP"
$ws.Range("N2").Value = ""

# -----------------------------------------------------------------
# 4) Row 3 (new order FINN18-002)
# -----------------------------------------------------------------
$ws.Range("A3").Value = 45439.083333333336
$ws.Range("B3").Value = "FINN18-002"
$ws.Range("C3").Value = "SIPL6153"
$ws.Range("D3").Value = "SIPL5688"
$ws.Range("E3").Value = "SIPL0102"
$ws.Range("F3").Value = "SIPL0103"
$ws.Range("G3").Value = "FINN TITLE"
$ws.Range("H3").Value = "Title"
$ws.Range("I3").Value = "Search & Typing"
$ws.Range("J3").Value = "Foreclosure information Report"
$ws.Range("K3").Value = "FL"
$ws.Range("L3").Value = "Clay"
$ws.Range("M3").Value = "WIP"
$ws.Range("N3").Value = ""

# -----------------------------------------------------------------
# 5) Column widths: autofit the columns whose contents changed size.
# -----------------------------------------------------------------
$ws.Range("C:C").EntireColumn.AutoFit() | Out-Null
$ws.Range("J:J").EntireColumn.AutoFit() | Out-Null
$ws.Range("N:N").EntireColumn.AutoFit() | Out-Null

# -----------------------------------------------------------------
# 6) Selection moved to H6 in the saved view.
# -----------------------------------------------------------------
$ws.Range("H6").Select() | Out-Null
